$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.293073296546936
$ws.Range("B1").Value = 2.89946460723877
$ws.Range("C1").Value = 5.174674034118652
$ws.Range("D1").Value = 1.859684824943542
$ws.Range("E1").Value = 1.011667728424072
